$d = $word.ActiveDocument

# 1. Rename the first paragraph's text: "Hello, world!" -> "Rebase 1".
$d.Content.Find.Execute("Hello, world!", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Rebase 1", 2)

# 2. Remove the "Change 1" and "Chanhe 2" paragraphs entirely, merging the
#    first paragraph directly with the trailing (bookmark) paragraph.
#    Merge paragraph marks one at a time (working from the end backwards so
#    earlier offsets stay valid), then strip the now-merged stray text.

# Merge paragraph 3 ("Chanhe 2") with paragraph 4 (the bookmark paragraph).
$p3end = $d.Paragraphs(3).Range.End
$d.Range($p3end - 1, $p3end).Delete()

# Merge paragraph 2 ("Change 1") with the combined paragraph above.
$p2end = $d.Paragraphs(2).Range.End
$d.Range($p2end - 1, $p2end).Delete()

# Delete the leftover "Change 1Chanhe 2" text that now lives in paragraph 2.
$p2 = $d.Paragraphs(2).Range
$p2Start = $p2.Start
$p2TextLen = $p2.Text.Length - 1   # exclude the trailing paragraph mark
$d.Range($p2Start, $p2Start + $p2TextLen).Delete()

# Merge paragraph 1 ("Rebase 1") with the now-empty paragraph 2, so the run
# and bookmark end up together in a single paragraph.
$p1end = $d.Paragraphs(1).Range.End
$d.Range($p1end - 1, $p1end).Delete()
